$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2328.8
$ws.Range("I40").Value = 2285.875
$ws.Range("K40").Value = 2285.875
$ws.Range("M40").Value = -2110.875
$ws.Range("H132").Value = 16951394
$ws.Range("I132").Value = 19232948
$ws.Range("J132").Value = 2695.2856
$ws.Range("K132").Value = 57698844
$ws.Range("L132").Value = 8085.8568
$ws.Range("M132").Value = -57696314
$ws.Range("N132").Value = -13145.8568
$ws.Range("H138").Value = 2593.795
$ws.Range("I138").Value = 1295.3636
$ws.Range("J138").Value = 4274.1177
$ws.Range("K138").Value = 3886.0908
$ws.Range("L138").Value = 12822.3531
$ws.Range("M138").Value = 1253.9092
$ws.Range("N138").Value = -23102.3531
$ws.Range("H141").Value = 2226.3684
$ws.Range("I141").Value = 1561.6765
$ws.Range("K141").Value = 4685.029500000001
$ws.Range("M141").Value = 494.9704999999994

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3260
$ws.Range("I74").Value = 2105.2727
$ws.Range("K74").Value = 2105.2727
$ws.Range("M74").Value = -1231.2727
$ws.Range("H77").Value = 3260
$ws.Range("I77").Value = 2105.2727
$ws.Range("K77").Value = 10526.3635
$ws.Range("M77").Value = -6158.363499999999
$ws.Range("H122").Value = 3960.6765
$ws.Range("I122").Value = 4122
$ws.Range("J122").Value = 3025
$ws.Range("K122").Value = 12366
$ws.Range("L122").Value = 9075
$ws.Range("M122").Value = -9916
$ws.Range("N122").Value = -13975

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 43479880
$ws.Range("I20").Value = 66668308
$ws.Range("J20").Value = 1573.875
$ws.Range("K20").Value = 66668308
$ws.Range("L20").Value = 1573.875
$ws.Range("M20").Value = -66668061
$ws.Range("N20").Value = -2067.875
$ws.Range("H94").Value = 111115304
$ws.Range("I94").Value = 2299.3333
$ws.Range("J94").Value = 166671800
$ws.Range("K94").Value = 2299.3333
$ws.Range("L94").Value = 166671800
$ws.Range("M94").Value = -1848.3333
$ws.Range("N94").Value = -166672702
$ws.Range("H134").Value = 4572.564
$ws.Range("I134").Value = 2315.3076
$ws.Range("K134").Value = 6945.9228
$ws.Range("M134").Value = -4410.9228

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H53").Value = 40061.332
$ws.Range("J53").Value = 40061.332
$ws.Range("L53").Value = 40061.332
$ws.Range("N53").Value = -41275.332
$ws.Range("H132").Value = 1474.579
$ws.Range("I132").Value = 1333.9375
$ws.Range("K132").Value = 4001.8125
$ws.Range("M132").Value = -1471.8125
$ws.Range("H134").Value = 1108.9734
$ws.Range("I134").Value = 1059.403
$ws.Range("J134").Value = 1524.125
$ws.Range("K134").Value = 3178.209
$ws.Range("L134").Value = 4572.375
$ws.Range("M134").Value = -643.2089999999998
$ws.Range("N134").Value = -9642.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 15151939
$ws.Range("I2").Value = 193.63637
$ws.Range("J2").Value = 30303684
$ws.Range("K2").Value = 1161.81822
$ws.Range("L2").Value = 181822104
$ws.Range("M2").Value = -1048.81822
$ws.Range("N2").Value = -181822330
$ws.Range("H5").Value = 1175
$ws.Range("I5").Value = 350
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 1050
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = -938
$ws.Range("N5").Value = -6224
$ws.Range("H38").Value = 391.9
$ws.Range("J38").Value = 539.5
$ws.Range("L38").Value = 1618.5
$ws.Range("N38").Value = -2312.5
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("H104").Value = 7899
$ws.Range("I104").Value = 7899
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 23697
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -21076
$ws.Range("N104").ClearContents()
$ws.Range("H105").Value = 2000
$ws.Range("J105").Value = 2000
$ws.Range("L105").Value = 6000
$ws.Range("N105").Value = -11242
$ws.Range("H106").Value = 24980
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H111").Value = 9999
$ws.Range("J111").Value = 9999
$ws.Range("L111").Value = 29997
$ws.Range("N111").Value = -36131
$ws.Range("H115").Value = 700055
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H118").Value = 1000000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 1000000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 3000000
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -3002486
$ws.Range("H120").Value = 13159.8
$ws.Range("I120").Value = 1799
$ws.Range("J120").Value = 16000
$ws.Range("K120").Value = 5397
$ws.Range("L120").Value = 48000
$ws.Range("M120").Value = -559
$ws.Range("N120").Value = -57676
$ws.Range("H135").Value = 1175
$ws.Range("I135").Value = 350
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 3150
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -615
$ws.Range("N135").Value = -23070

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 24019.5
$ws.Range("I70").Value = 33419.684
$ws.Range("J70").Value = 10280.77
$ws.Range("K70").Value = 33419.684
$ws.Range("L70").Value = 10280.77
$ws.Range("M70").Value = -33149.684
$ws.Range("N70").Value = -10820.77
$ws.Range("H73").Value = 24019.5
$ws.Range("I73").Value = 33419.684
$ws.Range("J73").Value = 10280.77
$ws.Range("K73").Value = 33419.684
$ws.Range("L73").Value = 10280.77
$ws.Range("M73").Value = -32483.684
$ws.Range("N73").Value = -12152.77
$ws.Range("H80").Value = 2964.3635
$ws.Range("I80").Value = 2461.8
$ws.Range("J80").Value = 3383.1667
$ws.Range("K80").Value = 2461.8
$ws.Range("L80").Value = 3383.1667
$ws.Range("M80").Value = -1463.8
$ws.Range("N80").Value = -5379.1667
$ws.Range("H83").Value = 2964.3635
$ws.Range("I83").Value = 2461.8
$ws.Range("J83").Value = 3383.1667
$ws.Range("K83").Value = 12309
$ws.Range("L83").Value = 16915.8335
$ws.Range("M83").Value = -7317
$ws.Range("N83").Value = -26899.8335
$ws.Range("H102").Value = 2565.3103
$ws.Range("I102").Value = 2510.44
$ws.Range("J102").Value = 2908.25
$ws.Range("K102").Value = 2510.44
$ws.Range("L102").Value = 2908.25
$ws.Range("M102").Value = -888.4400000000001
$ws.Range("N102").Value = -6152.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 865.7857
$ws.Range("I16").Value = 886.3077
$ws.Range("J16").Value = 599
$ws.Range("K16").Value = 886.3077
$ws.Range("L16").Value = 599
$ws.Range("M16").Value = -716.3077
$ws.Range("N16").Value = -939
$ws.Range("H46").Value = 790
$ws.Range("I46").Value = 832.8333
$ws.Range("J46").Value = 725.75
$ws.Range("K46").Value = 832.8333
$ws.Range("L46").Value = 725.75
$ws.Range("M46").Value = -644.8333
$ws.Range("N46").Value = -1101.75
$ws.Range("H93").Value = 4681.4585
$ws.Range("I93").Value = 4368.294
$ws.Range("J93").Value = 5442
$ws.Range("K93").Value = 4368.294
$ws.Range("L93").Value = 5442
$ws.Range("M93").Value = -3120.294
$ws.Range("N93").Value = -7938

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4658.884
$ws.Range("I122").Value = 4689.795
$ws.Range("J122").Value = 4357.5
$ws.Range("K122").Value = 14069.385
$ws.Range("L122").Value = 13072.5
$ws.Range("M122").Value = -11619.385
$ws.Range("N122").Value = -17972.5
$ws.Range("H132").Value = 1026.7179
$ws.Range("I132").Value = 701.8570999999999
$ws.Range("J132").Value = 1853.6364
$ws.Range("K132").Value = 2105.5713
$ws.Range("L132").Value = 5560.9092
$ws.Range("M132").Value = 424.4287000000004
$ws.Range("N132").Value = -10620.9092
